# "Controlo de progresso semana 9"
# Updates the weekly progress-control sheet: bumps the report date,
# marks several PT4/PT5 tasks as complete (or partially complete),
# fills in missing dates/owners, and fixes the view selection.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Report date (A2): 2024-ish serials, week bumped from 45645 to 45657 ---
$ws.Range("A2").Value2 = 45657

# --- PT3 block (rows 21-30): "Percentagem Execução Actual" (F) all reach 100% (1) ---
$ws.Range("F21").Value2 = 1
$ws.Range("F22").Value2 = 1
$ws.Range("F23").Value2 = 1
$ws.Range("F24").Value2 = 1
$ws.Range("F25").Value2 = 1

# F26 previously had a stray/inconsistent border style (s=25) vs. the rest of the
# block (s=19). Re-pull the formatting from a neighboring "correct" cell before
# updating the value so the style matches the rest of the column.
$ws.Range("F25").Copy()
$ws.Range("F26").PasteSpecial(-4122)
$ws.Range("F26").Value2 = 1

$ws.Range("F27").Value2 = 1
$ws.Range("F28").Value2 = 1
$ws.Range("F29").Value2 = 1
$ws.Range("F30").Value2 = 1

# --- PT4 block ---
# Row 32 (T4.1): finish date filled in, owner changed to "Todos", progress 0->0.3, actual 0.3->1
$ws.Range("C33").Copy()
$ws.Range("D32").PasteSpecial(-4122)
$ws.Range("D32").Value2 = 45649
$ws.Range("E32").Value2 = "Todos"
$ws.Range("F32").Value2 = 0.3
$ws.Range("G32").Value2 = 1

# Row 33 (T4.2): finish date filled in, owner changed to "Todos", progress 0->0.2, actual 0.2->1
$ws.Range("C33").Copy()
$ws.Range("D33").PasteSpecial(-4122)
$ws.Range("D33").Value2 = 45649
$ws.Range("E33").Value2 = "Todos"
$ws.Range("F33").Value2 = 0.2
$ws.Range("G33").Value2 = 1

# Row 34 (MP4.1): finish date filled in, owner set to "Rafael Fernandes", actual 0->1
$ws.Range("C34").Copy()
$ws.Range("D34").PasteSpecial(-4122)
$ws.Range("D34").Value2 = 45657
$ws.Range("E34").Value2 = "Rafael Fernandes"
$ws.Range("G34").Value2 = 1

# --- PT5 block ---
# Row 36 (T5.1): start date filled in, owner set to "Todos", actual 0->0.3
$ws.Range("C33").Copy()
$ws.Range("C36").PasteSpecial(-4122)
$ws.Range("C36").Value2 = 45657
$ws.Range("E36").Value2 = "Todos"
$ws.Range("G36").Value2 = 0.3

# Row 37 (T5.2): start date filled in
$ws.Range("C33").Copy()
$ws.Range("C37").PasteSpecial(-4122)
$ws.Range("C37").Value2 = 45657

# Row 38 (T5.3): start date filled in
$ws.Range("C33").Copy()
$ws.Range("C38").PasteSpecial(-4122)
$ws.Range("C38").Value2 = 45657

# Row 39 (T5.4): start+finish date filled in, owner "Rafael Fernandes", actual 0->1
$ws.Range("C33").Copy()
$ws.Range("C39").PasteSpecial(-4122)
$ws.Range("C39").Value2 = 45657
$ws.Range("C33").Copy()
$ws.Range("D39").PasteSpecial(-4122)
$ws.Range("D39").Value2 = 45657
$ws.Range("E39").Value2 = "Rafael Fernandes"
$ws.Range("G39").Value2 = 1

# Row 40 (T5.5): start+finish date filled in, owner "Todos", actual 0->1
$ws.Range("C33").Copy()
$ws.Range("C40").PasteSpecial(-4122)
$ws.Range("C40").Value2 = 45657
$ws.Range("C33").Copy()
$ws.Range("D40").PasteSpecial(-4122)
$ws.Range("D40").Value2 = 45657
$ws.Range("E40").Value2 = "Todos"
$ws.Range("G40").Value2 = 1

# Row 41 (T5.6): start+finish date filled in, owner "Todos", actual 0->1
$ws.Range("C33").Copy()
$ws.Range("C41").PasteSpecial(-4122)
$ws.Range("C41").Value2 = 45657
$ws.Range("C33").Copy()
$ws.Range("D41").PasteSpecial(-4122)
$ws.Range("D41").Value2 = 45657
$ws.Range("E41").Value2 = "Todos"
$ws.Range("G41").Value2 = 1

# Row 42 (T5.7): start date filled in
$ws.Range("C33").Copy()
$ws.Range("C42").PasteSpecial(-4122)
$ws.Range("C42").Value2 = 45657

# Row 43 (MP5.1): start date filled in
$ws.Range("C34").Copy()
$ws.Range("C43").PasteSpecial(-4122)
$ws.Range("C43").Value2 = 45657

# --- View state: selection moved to E41 ---
$ws.Range("E41").Select()
